$wb = $excel.ActiveWorkbook
$wsBets = $wb.Worksheets.Item("bets")
$wsResumen = $wb.Worksheets.Item("resumen")

# ---------------------------------------------------------------------------
# 1) Add two new bet rows (16 & 17) to the "bets" sheet.
#    Insert at row 16 so the row above's formatting (date style on B,
#    text style on J, percent style on M) carries down, then restore the
#    trailing spacer row (previously row 18, holding only a formatted B18)
#    back to its original position.
# ---------------------------------------------------------------------------
$wsBets.Rows("16:17").Insert(-4121, 0)
$wsBets.Range("B20").Cut($wsBets.Range("B18"))
$wsBets.Rows("19:20").Delete()

# Row 16 - AOE II Titans League Platinum / The Viper
$wsBets.Range("A16").Value = 15
$wsBets.Range("B16").Value = 45215
$wsBets.Range("C16").Value = 1
$wsBets.Range("D16").Formula = "=F15"
$wsBets.Range("E16").Value = 57
$wsBets.Range("F16").Formula = "=D16+E16"
$wsBets.Range("G16").Value = "ESPORTS"
$wsBets.Range("H16").Value = "AOE II TITANS LEAGUE PLATINUM"
$wsBets.Range("I16").Value = "THE VIPER"
$wsBets.Range("J16").Value = "GANA SERIE"
$wsBets.Range("K16").Value = 1
$wsBets.Range("L16").Value = 0
$wsBets.Range("M16").Formula = "=ROUND((F16/`$D`$2-1)*100, 3)"

# Row 17 - Qatar Masters / Magnus Carlsen
$wsBets.Range("A17").Value = 16
$wsBets.Range("B17").Value = 45216
$wsBets.Range("C17").Value = 1
$wsBets.Range("D17").Formula = "=F16"
$wsBets.Range("E17").Value = 84
$wsBets.Range("F17").Formula = "=D17+E17"
$wsBets.Range("G17").Value = "AJEDREZ"
$wsBets.Range("H17").Value = "QATAR MASTERS"
$wsBets.Range("I17").Value = "MAGNUS CARLSEN"
$wsBets.Range("J17").Value = "GANA O EMPATA"
$wsBets.Range("K17").Value = 1
$wsBets.Range("L17").Value = 0
$wsBets.Range("M17").Formula = "=ROUND((F17/`$D`$2-1)*100, 3)"

# Column H needs to grow to fit the new, longer category text.
$wsBets.Columns("H:H").ColumnWidth = 27.25

# ---------------------------------------------------------------------------
# 2) Selection / active-sheet bookkeeping: "bets" becomes the active tab
#    (it was "resumen" before), with H17 selected there; "resumen" keeps
#    its own prior selection but loses tab focus.
# ---------------------------------------------------------------------------
$wsResumen.Range("C6").Select()
$wsBets.Activate()
$wsBets.Range("H17").Select()

# resumen's array formulas reference bets!$M:$M through COUNTA/INDEX, so
# they automatically roll forward onto the new M17 once recalculated.
$excel.Calculate()
